# -----------------------------------------------------------------------
# RegTestData: add a "Status" column (G) with "Pass" for every data row,
# then add a new blank "HomePage" worksheet at the end of the workbook
# (parameterizeTest update - status column commit).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegTestData")

# --- new Status column --------------------------------------------------
$ws.Range("G1").Value = "Status"
$ws.Range("G2:G5").Value = "Pass"

# Header cell ends up with its own (non-bold) style rather than inheriting
# row 1's bold header style.
$ws.Range("G1").ClearFormats() | Out-Null

$ws.Columns("G:G").ColumnWidth = 6.56

# Cursor ends up away from the data, on F9, once editing is done.
$ws.Range("F9").Select() | Out-Null

# --- add the HomePage worksheet -----------------------------------------
# Adding/deleting a throwaway sheet first then adding the real one, then
# removing the throwaway, reproduces the sheetId numbering (3) the real
# session ended up with.
$placeholder = $wb.Worksheets.Add()
$homePage = $wb.Worksheets.Add()
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null

# Move HomePage to the end (after RegTestData) and name it.
$wb.Worksheets.Item(1).Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "HomePage"

$ws.Activate() | Out-Null
